$wb = $excel.ActiveWorkbook

# --- Sheets ---
$wsConv = $wb.Worksheets.Item("conversionAssets")   # sheet3.xml
$wsStor = $wb.Worksheets.Item("storageAssets")       # sheet4.xml

# ---------------------------------------------------------------------------
# conversionAssets: rename the gas burner asset and add the new DH gas burner
# ---------------------------------------------------------------------------

# Rename "House_gasburner" -> "House_gas_burner"
$wsConv.Range("B5").Value = "House_gas_burner"

# Add new default conversion asset row: DH_gas_burner_S
$wsConv.Range("A7").Value = 6
$wsConv.Range("B7").Value = "DH_gas_burner_S"
$wsConv.Range("C7").Value = "CONVERSION"
$wsConv.Range("D7").Value = "GAS_BURNER"
$wsConv.Range("E7").Value = 0
$wsConv.Range("F7").Value = 1000
$wsConv.Range("G7").Value = 0.95

# ---------------------------------------------------------------------------
# storageAssets: rename existing default heat-model / battery assets,
# backfill the new chargeCapacity_kWh column, and add a new default
# District Heating heat-buffer storage asset.
# ---------------------------------------------------------------------------

$wsStor.Range("B2").Value = "House_battery"
$wsStor.Range("B3").Value = "House_heatmodel_A"
$wsStor.Range("B4").Value = "House_heatmodel_B"
$wsStor.Range("B5").Value = "House_heatmodel_C"
$wsStor.Range("B6").Value = "House_heatmodel_D"
$wsStor.Range("B7").Value = "House_heatmodel_E"
$wsStor.Range("B8").Value = "House_heatmodel_F"
$wsStor.Range("B9").Value = "House_heatmodel_G"

# Backfill chargeCapacity_kWh (column K) for the heat-model rows
$wsStor.Range("K3").Value = 0
$wsStor.Range("K4").Value = 0
$wsStor.Range("K5").Value = 0
$wsStor.Range("K6").Value = 0
$wsStor.Range("K7").Value = 0
$wsStor.Range("K8").Value = 0
$wsStor.Range("K9").Value = 0

# Add new default storage asset row: District_Heating_heat_buffer_S
$wsStor.Range("A11").Value = 10
$wsStor.Range("B11").Value = "District_Heating_heat_buffer_S"
$wsStor.Range("C11").Value = "STORAGE"
$wsStor.Range("D11").Value = "STORAGE_HEAT"
$wsStor.Range("E11").Value = 0
$wsStor.Range("F11").Value = 100
$wsStor.Range("G11").Value = 0.5
$wsStor.Range("H11").Value = 0
$wsStor.Range("I11").Value = 90
$wsStor.Range("J11").Value = 1
$wsStor.Range("K11").Value = 0
$wsStor.Range("L11").Value = 10000

# ---------------------------------------------------------------------------
# View state: storageAssets becomes the active/selected sheet
# ---------------------------------------------------------------------------

$wsCons = $wb.Worksheets.Item("consumptionAssets")
$wsProd = $wb.Worksheets.Item("productionAssets")

$wsCons.Range("F11").Select()
$wsProd.Range("F8").Select()
$wsConv.Range("D9").Select()
$wsStor.Activate()
$wsStor.Range("I22").Select()
